# Refresh the "cryptos" price/volume table (GitHub Actions scheduled update).
# Updates Price (D) and Volume(1h) (E) for most rows, and for rows 31/32 the
# Toncoin / InjectiveProtocol entries swapped rank order, so their Coin (B),
# Link (C), Price (D) and Volume(1h) (E) all change together.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.312.73'
$ws.Range("E2").Value = '  -0.32%  '

$ws.Range("D3").Value = '2.600.70'
$ws.Range("E3").Value = '  +2.33%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '306.94'
$ws.Range("E5").Value = '  -0.01%  '

$ws.Range("D6").Value = '99.39'
$ws.Range("E6").Value = '  -4.10%  '

# D7/D32 are numeric-looking strings with a significant trailing zero
# ("0.600", "2.20"); a plain .Value assignment would auto-convert them to
# numbers and silently drop the trailing zero, so force text storage with a
# leading apostrophe and then strip the resulting quote-prefix style back to
# Normal so no stray number-format sticks to the cell.
$ws.Range("D7").Value = "'0.600"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.89%  '

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("E9").Value = '  +0.32%  '

$ws.Range("D10").Value = '39.28'
$ws.Range("E10").Value = '  +0.96%  '

$ws.Range("D11").Value = '54.13'
$ws.Range("E11").Value = '  -0.99%  '

$ws.Range("E12").Value = '  +1.45%  '

$ws.Range("D13").Value = '8.08'
$ws.Range("E13").Value = '  +1.38%  '

$ws.Range("D14").Value = '3.003.19'
$ws.Range("E14").Value = '  +2.42%  '

$ws.Range("E15").Value = '  +0.56%  '

$ws.Range("D16").Value = '2.603.91'
$ws.Range("E16").Value = '  +2.19%  '

$ws.Range("D17").Value = '0.915'
$ws.Range("E17").Value = '  +1.93%  '

$ws.Range("D18").Value = '14.88'
$ws.Range("E18").Value = '  -1.51%  '

$ws.Range("D19").Value = '46.378.73'
$ws.Range("E19").Value = '  -0.40%  '

$ws.Range("E20").Value = '  +1.19%  '

$ws.Range("D21").Value = '12.91'
$ws.Range("E21").Value = '  -8.08%  '

$ws.Range("D22").Value = '6.67'
$ws.Range("E22").Value = '  +1.03%  '

$ws.Range("D23").Value = '71.17'
$ws.Range("E23").Value = '  +1.88%  '

$ws.Range("D24").Value = '271.36'
$ws.Range("E24").Value = '  +6.62%  '

$ws.Range("E25").Value = '  +1.03%  '

$ws.Range("D26").Value = '2.16'
$ws.Range("E26").Value = '  +1.49%  '

$ws.Range("D27").Value = '29.15'
$ws.Range("E27").Value = '  +20.65%  '

$ws.Range("E28").Value = '  -0.01%  '

$ws.Range("E29").Value = '  -0.70%  '

$ws.Range("E30").Value = '  +0.49%  '

$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").Value = '38.44'
$ws.Range("E31").Value = '  -9.16%  '

$ws.Range("B32").Value = 'Toncoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D32").Value = "'2.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.91%  '

$ws.Range("E33").Value = '  +4.96%  '

$ws.Range("D34").Value = '3.63'
$ws.Range("E34").Value = '  -4.05%  '

$ws.Range("E35").Value = '  -2.03%  '

$ws.Range("E36").Value = '  +1.71%  '

$ws.Range("E37").Value = '  -1.26%  '

$ws.Range("D38").Value = '151.19'
$ws.Range("E38").Value = '  +0.29%  '

$ws.Range("E39").Value = '  +2.19%  '

$ws.Range("E40").Value = '  +1.20%  '

$ws.Range("D41").Value = '23.21'
$ws.Range("E41").Value = '  +31.80%  '

$ws.Range("D42").Value = '15.81'
$ws.Range("E42").Value = '  -4.82%  '

$ws.Range("D43").Value = '0.0329'
$ws.Range("E43").Value = '  +0.35%  '

$ws.Range("D44").Value = '3.57'
$ws.Range("E44").Value = '  +0.40%  '

$ws.Range("D45").Value = '4.05'
$ws.Range("E45").Value = '  -4.48%  '

$ws.Range("D46").Value = '2.115.11'
$ws.Range("E46").Value = '  +5.70%  '

$ws.Range("E47").Value = '  -0.04%  '

$ws.Range("D48").Value = '93.07'
$ws.Range("E48").Value = '  -0.90%  '

$ws.Range("D49").Value = '9.52'

$ws.Range("E50").Value = '  -5.80%  '

$ws.Range("D51").Value = '108.21'
$ws.Range("E51").Value = '  +1.15%  '
